$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D8").Value = 0.0862583333333
$ws.Range("E8").Value = 0.0547983823514
